$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.697.58"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").Value = "'3.186.74"
$ws.Range("E3").Value = "  +1.33%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'533.63"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("D6").Value = "'144.48"
$ws.Range("E6").Value = "  +3.63%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +2.13%  "
$ws.Range("D9").Value = "'7.31"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("E10").Value = "  +2.31%  "
$ws.Range("D11").Value = "'0.428"
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("D12").Value = "'3.736.77"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("E13").Value = "  -1.41%  "
$ws.Range("D14").Value = "'25.93"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").Value = "'59.791.41"
$ws.Range("E16").Value = "  +2.10%  "
$ws.Range("D17").Value = "'3.224.41"
$ws.Range("E17").Value = "  +2.05%  "
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("D20").Value = "'8.17"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").Value = "'366.00"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("D24").Value = "'69.52"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("E25").Value = "  +10.64%  "
$ws.Range("D26").Value = "'0.168"
$ws.Range("E26").Value = "  +1.13%  "
$ws.Range("D27").Value = "'0.997"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").Value = "'0.0₃0877"
$ws.Range("E28").Value = "  +2.74%  "
$ws.Range("D29").Value = "'22.23"
$ws.Range("E29").Value = "  +1.76%  "
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("D32").Value = "'5.27"
$ws.Range("E32").Value = "  +1.86%  "
$ws.Range("E33").Value = "  +2.40%  "
$ws.Range("D34").Value = "'6.54"
$ws.Range("E34").Value = "  +4.59%  "
$ws.Range("D35").Value = "'156.43"
$ws.Range("E35").Value = "  -1.63%  "
$ws.Range("E36").Value = "  -1.57%  "
$ws.Range("D37").Value = "'2.784.82"
$ws.Range("E37").Value = "  +5.33%  "
$ws.Range("D38").Value = "'25.83"
$ws.Range("E38").Value = "  +2.96%  "
$ws.Range("D39").Value = "'0.0697"
$ws.Range("E39").Value = "  +2.26%  "
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'39.42"
$ws.Range("E42").Value = "  +2.38%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0291"
$ws.Range("E43").Value = "  +2.42%  "
$ws.Range("E44").Value = "  +0.55%  "
$ws.Range("D45").Value = "'0.105"
$ws.Range("E45").Value = "  +3.43%  "
$ws.Range("D46").Value = "'3.227.80"
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("D48").Value = "'0.805"
$ws.Range("E48").Value = "  +7.12%  "
$ws.Range("E49").Value = "  -0.89%  "
$ws.Range("D50").Value = "'20.42"
$ws.Range("E50").Value = "  +2.06%  "
$ws.Range("E51").Value = "  +0.03%  "
